$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cryptocurrency price/volume/name/link cells.
# A leading apostrophe forces Excel to store numeric-looking text (e.g. '96.166.98')
# as a text string rather than auto-converting it to a number.

$ws.Range("D2").Value = '''96.166.98'
$ws.Range("E2").Value = '''  +4.41%  '
$ws.Range("D3").Value = '''3.654.46'
$ws.Range("E3").Value = '''  +9.64%  '
$ws.Range("E4").Value = '''  +0.04%  '
$ws.Range("D5").Value = '''240.15'
$ws.Range("E5").Value = '''  +3.88%  '
$ws.Range("D6").Value = '''643.98'
$ws.Range("E6").Value = '''  +4.52%  '
$ws.Range("D7").Value = '''1.48'
$ws.Range("E7").Value = '''  +5.49%  '
$ws.Range("E8").Value = '''  +3.79%  '
$ws.Range("E9").Value = '''  -0.14%  '
$ws.Range("E10").Value = '''  +5.76%  '
$ws.Range("D11").Value = '''3.656.80'
$ws.Range("E11").Value = '''  +9.76%  '
$ws.Range("D12").Value = '''43.40'
$ws.Range("E12").Value = '''  +1.11%  '
$ws.Range("E13").Value = '''  +3.35%  '
$ws.Range("D14").Value = '''6.35'
$ws.Range("E14").Value = '''  +4.93%  '
$ws.Range("D15").Value = '''4.344.00'
$ws.Range("E15").Value = '''  +9.67%  '
$ws.Range("D16").Value = '''96.114.06'
$ws.Range("E16").Value = '''  +4.51%  '
$ws.Range("E17").Value = '''  +5.00%  '
$ws.Range("B18").Value = '''Uniswap'
$ws.Range("C18").Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '''13.54'
$ws.Range("E18").Value = '''  +24.29%  '
$ws.Range("B19").Value = '''WrappedEther'
$ws.Range("C19").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '''3.645.79'
$ws.Range("E19").Value = '''  +9.28%  '
$ws.Range("E20").Value = '''  -1.06%  '
$ws.Range("D21").Value = '''18.31'
$ws.Range("E21").Value = '''  +5.68%  '
$ws.Range("B22").Value = '''BitcoinCash'
$ws.Range("C22").Value = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '''517.73'
$ws.Range("E22").Value = '''  +4.44%  '
$ws.Range("B23").Value = '''Stellar'
$ws.Range("C23").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D23").Value = '''0.485'
$ws.Range("E23").Value = '''  +8.55%  '
$ws.Range("E24").Value = '''  +0.06%  '
$ws.Range("E25").Value = '''  +7.92%  '
$ws.Range("D26").Value = '''6.70'
$ws.Range("E26").Value = '''  +7.11%  '
$ws.Range("D27").Value = '''97.35'
$ws.Range("E27").Value = '''  +8.43%  '
$ws.Range("D28").Value = '''12.55'
$ws.Range("E28").Value = '''  +5.78%  '
$ws.Range("D29").Value = '''3.14'
$ws.Range("E29").Value = '''  +20.45%  '
$ws.Range("D30").Value = '''11.63'
$ws.Range("E30").Value = '''  +4.99%  '
$ws.Range("D31").Value = '''0.142'
$ws.Range("E31").Value = '''  +2.43%  '
$ws.Range("E32").Value = '''  -0.09%  '
$ws.Range("E33").Value = '''  +5.76%  '
$ws.Range("E34").Value = '''  +1.60%  '
$ws.Range("D35").Value = '''32.05'
$ws.Range("E35").Value = '''  +13.36%  '
$ws.Range("D36").Value = '''0.577'
$ws.Range("E36").Value = '''  +9.24%  '
$ws.Range("D37").Value = '''566.70'
$ws.Range("E37").Value = '''  -0.41%  '
$ws.Range("D38").Value = '''7.84'
$ws.Range("E38").Value = '''  +6.85%  '
$ws.Range("E39").Value = '''  +8.00%  '
$ws.Range("D40").Value = '''0.950'
$ws.Range("E40").Value = '''  +9.55%  '
$ws.Range("E41").Value = '''  +3.24%  '
$ws.Range("E42").Value = '''  -0.08%  '
$ws.Range("E43").Value = '''  +4.02%  '
$ws.Range("B44").Value = '''Filecoin'
$ws.Range("C44").Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '''5.77'
$ws.Range("E44").Value = '''  +6.49%  '
$ws.Range("B45").Value = '''ImmutableX'
$ws.Range("C45").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").Value = '''1.73'
$ws.Range("E45").Value = '''  +3.73%  '
$ws.Range("E46").Value = '''  +0.40%  '
$ws.Range("E47").Value = '''  +5.36%  '
$ws.Range("B48").Value = '''OKB'
$ws.Range("C48").Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").Value = '''54.14'
$ws.Range("E48").Value = '''  +5.04%  '
$ws.Range("B49").Value = '''MantraDAO'
$ws.Range("C49").Value = '''https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D49").Value = '''3.49'
$ws.Range("E49").Value = '''  -3.05%  '
$ws.Range("D50").Value = '''8.20'
$ws.Range("E50").Value = '''  +2.89%  '
$ws.Range("E51").Value = '''  +4.01%  '

# Clear the "quote prefix" text-override styling that typing a leading
# apostrophe adds, so number formats/styles stay exactly as they were.
$ws.Range("B2:E51").Style = "Normal"
